$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.212.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.305.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.68%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.505'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.76%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.516'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.28'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.72%  '
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.86'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +12.36%  '
$ws.Range("E13").Value = '  +1.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.664.67'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.264.05'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("E17").Value = '  +0.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.058.71'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +10.75%  '
$ws.Range("E20").Value = '  +3.74%  '
$ws.Range("E21").Value = '  +1.25%  '
$ws.Range("E22").Value = '  +1.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.85'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("E24").Value = '  +13.95%  '
$ws.Range("E25").Value = '  +0.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.46'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.36'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.93'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '167.47'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.88%  '
$ws.Range("E32").Value = '  +0.07%  '
$ws.Range("E33").Value = '  +1.97%  '
$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.72'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.35%  '
$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.82'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.71%  '
$ws.Range("E36").Value = '  +0.74%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0698'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.64%  '
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("E39").Value = '  +3.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.101'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.32%  '
$ws.Range("E41").Value = '  +1.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.38'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.998.53'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.03%  '
$ws.Range("E44").Value = '  +3.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.17'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.70'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("E47").Value = '  +2.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.80'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.38%  '
$ws.Range("E49").Value = '  +5.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.531.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '70.95'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.74%  '
